$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect, make the edits, then restore protection.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer string (A9).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."
# Setting a multi-line value can trigger an automatic row-height resize;
# auto-fit the row back so it doesn't pick up a stray explicit height.
$ws.Rows(9).AutoFit()

# Updated Weight / Percent Change figures for the holdings rows.
$ws.Range("D2").Value = 0.2514910377617529
$ws.Range("E2").Value = 0.01091599430469881

$ws.Range("D3").Value = 0.2509521607209376
$ws.Range("E3").Value = 0.007958287596048219

$ws.Range("D4").Value = 0.2463590023950406
$ws.Range("E4").Value = 0.004369780539910773

$ws.Range("D5").Value = 0.251197799122269
$ws.Range("E5").Value = -0.009015971148892299

$ws.Range("E6").Value = 0.003554166868703135

# Restore sheet protection.
$ws.Protect()
